$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update QUERYSTRING (column G) and VALIDATIONS (column J) cells to replace
# the "mohana"/"moha"/"phani" test data with "Project"/"Proj"/"Neon2"/"Neon1".

$ws.Range("G2").Value = "?query=Project&size=1"
$ws.Range("J2").Value = "status=200||hits.hits._source.firstName=Project//hits.hits._source.lastName=Project"

$ws.Range("G6").Value = "?query=Project&size=2"
$ws.Range("J6").Value = "status=200||hits.hits._source.firstName=Project//hits.hits._source.lastName=Project"

$ws.Range("G7").Value = "?query=Project+Neon2"
$ws.Range("J7").Value = "status=200||hits.hits._source.firstName=Project//hits.hist._source.lastName=Project//hits.hits._source.firstName=Neon2//hits.hits._source.lastName=Neon2"

$ws.Range("G8").Value = "?query=Proj*"
$ws.Range("J8").Value = "status=200||hits.hits._source.firstName=Project//hits.hits._source.lastName=Project"

$ws.Range("G9").Value = "?query=Proj* -Neon2"
$ws.Range("J9").Value = "status=200||hits.hits._source.firstName=Project//hits.hits._source.lastName=Project"

$ws.Range("G10").Value = "?query=Proj* -Neon2 -Neon1"
$ws.Range("J10").Value = "status=200"

$ws.Range("G11").Value = "?query=Project&size=2&fields=category"

$wb.Save()
